$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the rich-text header of E1: drop the trailing line break after the
# closing parenthesis in "Giá trị tối thiểu của hóa đơn (*)\n" -> "...(*)"
# while preserving the per-character formatting (bold red "*", bold ")").
$cell = $ws.Range("E1")
$cell.Value = "Giá trị tối thiểu của hóa đơn (*)"
$star = $cell.Characters(32, 1)
$star.Font.Bold = $true
$star.Font.Color = 255
$closingParen = $cell.Characters(33, 1)
$closingParen.Font.Bold = $true

# --- Update the window/view state: scroll so column B is at the left edge
# and move the active selection to E2 (previously F1 / J13).
$win = $excel.ActiveWindow
$win.ScrollColumn = 2
$win.ScrollRow = 1
$win.Left = -110
$win.Top = -110
$win.Width = 19420
$win.Height = 10300

$ws.Range("E2").Select() | Out-Null
